# Created UI for input panel config in Time mode.
# Adds 6 new rows (22-27) to the "Translation" sheet describing the new
# input-panel text ids used by the Time mode UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 22: Detect slope label
$ws.Range("B22").Value = "SingleUseId19"
$ws.Range("C22").Value = "Default"
$ws.Range("D22").Value = "Left"
$ws.Range("E22").Value = "LTR"
$ws.Range("F22").Value = "SLOPE"

# Row 23: Detect threshold label
$ws.Range("B23").Value = "SingleUseId20"
$ws.Range("C23").Value = "Default"
$ws.Range("D23").Value = "Left"
$ws.Range("E23").Value = "LTR"
$ws.Range("F23").Value = "THRESHOLD"

# Row 24: Detect mode option
$ws.Range("B24").Value = "SingleUseId21"
$ws.Range("C24").Value = "Default"
$ws.Range("D24").Value = "Center"
$ws.Range("E24").Value = "LTR"
$ws.Range("F24").Value = "Detect"

# Row 25: value placeholder
$ws.Range("B25").Value = "SingleUseId22"
$ws.Range("C25").Value = "Default"
$ws.Range("D25").Value = "Center"
$ws.Range("E25").Value = "LTR"
$ws.Range("F25").Value = "<value>"

# Row 26: Manual mode option
$ws.Range("B26").Value = "SingleUseId23"
$ws.Range("C26").Value = "Default"
$ws.Range("D26").Value = "Left"
$ws.Range("E26").Value = "LTR"
$ws.Range("F26").Value = "Manual"

# Row 27: Defined label
$ws.Range("B27").Value = "SingleUseId24"
$ws.Range("C27").Value = "Default"
$ws.Range("D27").Value = "Left"
$ws.Range("E27").Value = "LTR"
$ws.Range("F27").Value = "Defined"
